{"js": "// Update each three-digit-by-one-digit multiplication answer cell in the\n// table to the newly generated problem/answer pair. Each \"old\" value is\n// unique in the document, so a targeted search + replace is safe and\n// leaves every other run (fonts, sizes, paragraph props, empty rows)\n// untouched.\nconst replacements = [\n  [\"264\u00d72=528\", \"162\u00d78=1296\"],\n  [\"978\u00d78=7824\", \"967\u00d78=7736\"],\n  [\"398\u00d78=3184\", \"521\u00d78=4168\"],\n  [\"756\u00d76=4536\", \"652\u00d78=5216\"],\n  [\"385\u00d78=3080\", \"116\u00d72=232\"],\n  [\"897\u00d76=5382\", \"724\u00d77=5068\"],\n  [\"968\u00d74=3872\", \"577\u00d76=3462\"],\n  [\"198\u00d79=1782\", \"726\u00d76=4356\"],\n  [\"679\u00d74=2716\", \"113\u00d72=226\"],\n  [\"566\u00d72=1132\", \"538\u00d77=3766\"],\n  [\"529\u00d79=4761\", \"885\u00d78=7080\"],\n  [\"434\u00d78=3472\", \"554\u00d78=4432\"],\n  [\"423\u00d74=1692\", \"207\u00d72=414\"],\n  [\"183\u00d79=1647\", \"433\u00d73=1299\"],\n  [\"448\u00d75=2240\", \"140\u00d72=280\"],\n  [\"102\u00d79=918\", \"527\u00d75=2635\"],\n  [\"143\u00d76=858\", \"641\u00d79=5769\"],\n  [\"272\u00d77=1904\", \"738\u00d75=3690\"],\n  [\"664\u00d75=3320\", \"806\u00d74=3224\"],\n  [\"299\u00d75=1495\", \"264\u00d77=1848\"],\n  [\"926\u00d76=5556\", \"267\u00d73=801\"],\n  [\"138\u00d74=552\", \"601\u00d75=3005\"],\n  [\"163\u00d73=489\", \"887\u00d73=2661\"],\n  [\"423\u00d77=2961\", \"287\u00d75=1435\"],\n  [\"813\u00d77=5691\", \"439\u00d72=878\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication answer cell in the\n# table to the newly generated problem/answer pair. Each \"old\" value is\n# unique in the document, so a targeted Find/Replace is safe and leaves\n# every other run (fonts, sizes, paragraph props, empty rows) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"264\u00d72=528\", \"162\u00d78=1296\"),\n  @(\"978\u00d78=7824\", \"967\u00d78=7736\"),\n  @(\"398\u00d78=3184\", \"521\u00d78=4168\"),\n  @(\"756\u00d76=4536\", \"652\u00d78=5216\"),\n  @(\"385\u00d78=3080\", \"116\u00d72=232\"),\n  @(\"897\u00d76=5382\", \"724\u00d77=5068\"),\n  @(\"968\u00d74=3872\", \"577\u00d76=3462\"),\n  @(\"198\u00d79=1782\", \"726\u00d76=4356\"),\n  @(\"679\u00d74=2716\", \"113\u00d72=226\"),\n  @(\"566\u00d72=1132\", \"538\u00d77=3766\"),\n  @(\"529\u00d79=4761\", \"885\u00d78=7080\"),\n  @(\"434\u00d78=3472\", \"554\u00d78=4432\"),\n  @(\"423\u00d74=1692\", \"207\u00d72=414\"),\n  @(\"183\u00d79=1647\", \"433\u00d73=1299\"),\n  @(\"448\u00d75=2240\", \"140\u00d72=280\"),\n  @(\"102\u00d79=918\", \"527\u00d75=2635\"),\n  @(\"143\u00d76=858\", \"641\u00d79=5769\"),\n  @(\"272\u00d77=1904\", \"738\u00d75=3690\"),\n  @(\"664\u00d75=3320\", \"806\u00d74=3224\"),\n  @(\"299\u00d75=1495\", \"264\u00d77=1848\"),\n  @(\"926\u00d76=5556\", \"267\u00d73=801\"),\n  @(\"138\u00d74=552\", \"601\u00d75=3005\"),\n  @(\"163\u00d73=489\", \"887\u00d73=2661\"),\n  @(\"423\u00d77=2961\", \"287\u00d75=1435\"),\n  @(\"813\u00d77=5691\", \"439\u00d72=878\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute(\n    $oldText,  # FindText\n    $true,     # MatchCase\n    $false,    # MatchWholeWord\n    $false,    # MatchWildcards\n    $false,    # MatchSoundsLike\n    $false,    # MatchAllWordForms\n    $true,     # Forward\n    1,         # Wrap = wdFindContinue\n    $false,    # Format\n    $newText,  # ReplaceWith\n    2          # Replace = wdReplaceAll\n  )\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
